$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new columns (this shifts existing cell content + column-width
# metadata to the right, preserving the original column widths exactly) ---

# 3 new columns before old column G ("Pause"): will become "Dienstbeginn Ort",
# "Wartezeit 1", "Wartezeit 2" slots (G,H,I) -- "Pause" itself ends up at J (not
# final position yet, more inserts below).
$ws.Range("G1:I1").EntireColumn.Insert()

# 1 new column right after "Pause" (now column J) for "Pause Gesamt"
$ws.Columns.Item(11).Insert()

# 1 new column right after "Dienstende" (now column L) for "Dienstende Ort"
$ws.Columns.Item(13).Insert()

# --- Column widths ---
# Column D width changes (and loses the bestFit/auto-fit flag). The engine
# quantizes stored column widths to a 1/6-character-unit pixel grid, so the
# ColumnWidth input is pre-compensated to land on the closest achievable
# stored width to the target (62.7265625).
$ws.Columns.Item(4).ColumnWidth = 61.833333333333336

# New columns get their header widths (again pre-compensated for the 1/6 grid)
$ws.Columns.Item(7).ColumnWidth = 14.333333333333334
$ws.Columns.Item(8).ColumnWidth = 14.333333333333334
$ws.Columns.Item(9).ColumnWidth = 14.333333333333334
$ws.Columns.Item(11).ColumnWidth = 12.0
$ws.Columns.Item(13).ColumnWidth = 13.0

# --- Row 5 header text (full rewrite to the final, target layout) ---
$ws.Range("A5").Value = "Lfn"
$ws.Range("B5").Value = "Schicht Zug"
$ws.Range("C5").Value = "Datum"
$ws.Range("D5").Value = "Gastfahrt vor Dienstbeginn"
$ws.Range("E5").Value = "Gesamt"
$ws.Range("F5").Value = "Dienstbeginn"
$ws.Range("G5").Value = "Dienstbeginn Ort"
$ws.Range("H5").Value = "Pause"
$ws.Range("I5").Value = "Wartezeit 1"
$ws.Range("J5").Value = "Wartezeit 2"
$ws.Range("K5").Value = "Pause Gesamt"
$ws.Range("L5").Value = "Dienstende"
$ws.Range("M5").Value = "Dienstende Ort"
$ws.Range("N5").Value = "Abfahrt / Ankunft"
$ws.Range("O5").Value = "Gesamt"
$ws.Range("P5").Value = "Gastfahrt nach Dienstende"
$ws.Range("Q5").Value = "Gesamt"
$ws.Range("R5").Value = "Gesamt Dienststunden"
$ws.Range("S5").Value = "Bemerkungen"

# --- Selection moves to D5 ---
$ws.Range("D5").Select()
